$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# The (hidden) "_GoBack" bookmark currently sits right at the end of the
# "2 Change commented line ..." paragraph, at the tail of the 0.6 change-log
# block. The new 0.7 entry being merged in needs that bookmark moved to the
# end of its own last line, so drop the existing one first - it gets
# re-created below once the new content is in place.
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# ---------------------------------------------------------------------------
# Locate "2 Change commented line to properly run the SQLPSX.psm1." and the
# paragraph it lives in, independent of absolute paragraph indices.
# ---------------------------------------------------------------------------
$findRange = $d.Content
$found = $findRange.Find.Execute(
    "2 Change commented line to properly run the SQLPSX.psm1.",
    $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

$anchorStart = $findRange.Start
$anchorEnd = $findRange.End

$pAnchor = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Start -le $anchorStart -and $cand.Range.End -ge $anchorEnd) {
        $pAnchor = $cand
        break
    }
}

# The two paragraphs right after it are the blank, bold-formatted
# placeholder paragraphs left ready for the next change-log entry.
$pBlank1 = $pAnchor.Next()
$pBlank2 = $pBlank1.Next()

# Insert a brand-new paragraph right after the second placeholder - this
# becomes the "1. Found bug ..." line, and is where "_GoBack" ends up.
$pBlank2.Range.InsertParagraphAfter()
$pNew = $pBlank2.Next()

# Turn the (previously blank) second placeholder into the change-log header.
$pBlank2.Range.Text = "Change log - 0.7 - 05/10/2010 17:20 - Max Trinidad"
$pBlank2.Range.Bold = 1

# Fill in the bug-fix description on the freshly inserted paragraph.
$pNew.Range.Text = "1. Found bug missing path to Windows\system32 PowerShell modules."
$pNew.Range.Bold = 1

# ---------------------------------------------------------------------------
# Re-create "_GoBack" as a collapsed (zero-length) bookmark right after the
# new text, before the paragraph mark - matching its original shape. A
# temporary character anchors the bookmark to a non-empty range, then gets
# removed again, leaving the bookmark collapsed in place.
# ---------------------------------------------------------------------------
$endPos = $pNew.Range.End - 1
$d.Range($endPos, $endPos).InsertAfter("~")

$pNew = $pBlank2.Next()
$bmEnd = $pNew.Range.End - 1
$bmStart = $bmEnd - 1
$d.Bookmarks.Add("_GoBack", $d.Range($bmStart, $bmEnd))

$pNew = $pBlank2.Next()
$delStart = $pNew.Range.End - 2
$delEnd = $pNew.Range.End - 1
$d.Range($delStart, $delEnd).Text = ""

Write-Host "Change log 0.7 entry merged."
